$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1912.6842
$ws.Range("I40").Value = 2020.0834
$ws.Range("J40").Value = 1728.5714
$ws.Range("K40").Value = 2020.0834
$ws.Range("L40").Value = 1728.5714
$ws.Range("M40").Value = -1845.0834
$ws.Range("N40").Value = -2078.5714

$ws.Range("H76").Value = 4462.5625
$ws.Range("I76").Value = 3372.7273
$ws.Range("J76").Value = 6860.2
$ws.Range("K76").Value = 3372.7273
$ws.Range("L76").Value = 6860.2
$ws.Range("M76").Value = -3057.7273
$ws.Range("N76").Value = -7490.2

$ws.Range("H79").Value = 4462.5625
$ws.Range("I79").Value = 3372.7273
$ws.Range("J79").Value = 6860.2
$ws.Range("K79").Value = 3372.7273
$ws.Range("L79").Value = 6860.2
$ws.Range("M79").Value = -2280.7273
$ws.Range("N79").Value = -9044.200000000001

$ws.Range("H138").Value = 4552.919
$ws.Range("I138").Value = 1878.625
$ws.Range("J138").Value = 6590.476
$ws.Range("K138").Value = 5635.875
$ws.Range("L138").Value = 19771.428
$ws.Range("M138").Value = -495.875
$ws.Range("N138").Value = -30051.428

$ws.Range("H141").Value = 1963.6451
$ws.Range("I141").Value = 2030.138
$ws.Range("J141").Value = 999.5
$ws.Range("K141").Value = 6090.414
$ws.Range("L141").Value = 2998.5
$ws.Range("M141").Value = -910.4139999999998
$ws.Range("N141").Value = -13358.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1875.2307
$ws.Range("I86").Value = 1644.7059
$ws.Range("J86").Value = 2310.6667
$ws.Range("K86").Value = 1644.7059
$ws.Range("L86").Value = 2310.6667
$ws.Range("M86").Value = -521.7058999999999
$ws.Range("N86").Value = -4556.6667

$ws.Range("H89").Value = 1875.2307
$ws.Range("I89").Value = 1644.7059
$ws.Range("J89").Value = 2310.6667
$ws.Range("K89").Value = 8223.529500000001
$ws.Range("L89").Value = 11553.3335
$ws.Range("M89").Value = -2607.529500000001
$ws.Range("N89").Value = -22785.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27781684
$ws.Range("I31").Value = 45456028
$ws.Range("J31").Value = 7716.2856
$ws.Range("K31").Value = 45456028
$ws.Range("L31").Value = 7716.2856
$ws.Range("M31").Value = -45455733
$ws.Range("N31").Value = -8306.285599999999

$ws.Range("H34").Value = 27781684
$ws.Range("I34").Value = 45456028
$ws.Range("J34").Value = 7716.2856
$ws.Range("K34").Value = 45456028
$ws.Range("L34").Value = 7716.2856
$ws.Range("M34").Value = -45455826
$ws.Range("N34").Value = -8120.2856

$ws.Range("H54").Value = 37333.332
$ws.Range("J54").Value = 37333.332
$ws.Range("L54").Value = 37333.332
$ws.Range("N54").Value = -38649.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3875.6667
$ws.Range("I110").Value = 3875.6667
$ws.Range("K110").Value = 11627.0001
$ws.Range("M110").Value = -7537.000100000001

$ws.Range("H111").Value = 3601.889
$ws.Range("I111").Value = 2871.1667
$ws.Range("J111").Value = 5063.3335
$ws.Range("K111").Value = 8613.500100000001
$ws.Range("L111").Value = 15190.0005
$ws.Range("M111").Value = -5546.500100000001
$ws.Range("N111").Value = -21324.0005

$ws.Range("H114").Value = 1770
$ws.Range("I114").Value = 855.55554
$ws.Range("J114").Value = 10000
$ws.Range("K114").Value = 2566.66662
$ws.Range("L114").Value = 30000
$ws.Range("M114").Value = 687.33338
$ws.Range("N114").Value = -36508

$ws.Range("H115").Value = 1558.9524
$ws.Range("I115").Value = 1287.8918
$ws.Range("J115").Value = 3564.8
$ws.Range("K115").Value = 3863.6754
$ws.Range("L115").Value = 10694.4
$ws.Range("M115").Value = -2688.6754
$ws.Range("N115").Value = -13044.4

$ws.Range("H116").Value = 917.2
$ws.Range("I116").Value = 396.5
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1189.5
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = 2252.5
$ws.Range("N116").Value = -15884

$ws.Range("H117").Value = 1780
$ws.Range("I117").Value = 373.33334
$ws.Range("J117").Value = 6000
$ws.Range("K117").Value = 1120.00002
$ws.Range("L117").Value = 18000
$ws.Range("M117").Value = 2321.99998
$ws.Range("N117").Value = -24884

$ws.Range("H118").Value = 2140490.5
$ws.Range("I118").Value = 3375
$ws.Range("J118").Value = 3090319.5
$ws.Range("K118").Value = 10125
$ws.Range("L118").Value = 9270958.5
$ws.Range("M118").Value = -8882
$ws.Range("N118").Value = -9273444.5

$ws.Range("H119").Value = 2980
$ws.Range("I119").Value = 2475
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 7425
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = -2587
$ws.Range("N119").Value = -24676

$ws.Range("H120").Value = 18857.143
$ws.Range("I120").Value = 14000
$ws.Range("J120").Value = 19666.666
$ws.Range("K120").Value = 42000
$ws.Range("L120").Value = 58999.99800000001
$ws.Range("M120").Value = -37162
$ws.Range("N120").Value = -68675.99800000001

$ws.Range("H121").Value = 38466970
$ws.Range("I121").Value = 2500
$ws.Range("J121").Value = 40005548
$ws.Range("K121").Value = 7500
$ws.Range("L121").Value = 120016644
$ws.Range("M121").Value = -6190
$ws.Range("N121").Value = -120019264

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 10000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -10224

$ws.Range("H44").Value = 98000
$ws.Range("J44").Value = 98000
$ws.Range("L44").Value = 98000
$ws.Range("N44").Value = -99192

$ws.Range("H70").Value = 5579
$ws.Range("I70").Value = 5153.5
$ws.Range("J70").Value = 6004.5
$ws.Range("K70").Value = 5153.5
$ws.Range("L70").Value = 6004.5
$ws.Range("M70").Value = -4883.5
$ws.Range("N70").Value = -6544.5

$ws.Range("H73").Value = 5579
$ws.Range("I73").Value = 5153.5
$ws.Range("J73").Value = 6004.5
$ws.Range("K73").Value = 5153.5
$ws.Range("L73").Value = 6004.5
$ws.Range("M73").Value = -4217.5
$ws.Range("N73").Value = -7876.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13000
$ws.Range("J54").Value = 13000
$ws.Range("L54").Value = 13000
$ws.Range("N54").Value = -14040
